# Applies the "cambios en el excel de horarios parte 2" edit:
#  - Fills in previously-blank rows 34-37 on the "Sprint2" sheet
#  - Adds 4 new rows (38-41) with the same formatting as row 37
#  - Updates the saved selection on the "Sprint2" sheet to E26

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sprint2")

# --- Row 34 ---
$ws2.Range("C34").Value = "Miguel Ángel"
$ws2.Range("D34").Value = "miguelangelmateos8"
$ws2.Range("E34").Value = "Cambis en la base de dades i en els models"
$ws2.Range("F34").Value = "Leonard Craciun"
$ws2.Range("G34").Value = 3
$ws2.Range("H34").Value = 1
$ws2.Range("I34").Value = 3

# --- Row 35 ---
$ws2.Range("C35").Value = "Miguel Ángel"
$ws2.Range("D35").Value = "miguelangelmateos8"
$ws2.Range("E35").Value = "Vista de localitzacions"
$ws2.Range("F35").Value = "Leonard Craciun"
$ws2.Range("G35").Value = 1
$ws2.Range("H35").Value = 1
$ws2.Range("I35").Value = 1

# --- Row 36 ---
$ws2.Range("C36").Value = "Miguel Ángel"
$ws2.Range("D36").Value = "miguelangelmateos8"
$ws2.Range("E36").Value = "Logica de localitzacions"
$ws2.Range("F36").Value = "Leonard Craciun"
$ws2.Range("G36").Value = 1
$ws2.Range("H36").Value = 1
$ws2.Range("I36").Value = 1.5

# --- Row 37 ---
$ws2.Range("C37").Value = "Miguel Ángel"
$ws2.Range("D37").Value = "miguelangelmateos8"
$ws2.Range("E37").Value = "Crear localitzacions"
$ws2.Range("F37").Value = "Leonard Craciun"
$ws2.Range("G37").Value = 1.5
$ws2.Range("H37").Value = 1
$ws2.Range("I37").Value = 1

# --- New rows 38-41: copy the formatting from row 37 first ---
$ws2.Range("B37:I37").Copy()
$ws2.Range("B38:I41").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 38 ---
$ws2.Range("B38").Value = 2
$ws2.Range("C38").Value = "Miguel Ángel"
$ws2.Range("D38").Value = "miguelangelmateos8"
$ws2.Range("E38").Value = "Vista exposicions"
$ws2.Range("F38").Value = "Leonard Craciun"
$ws2.Range("G38").Value = 2
$ws2.Range("H38").Value = 1
$ws2.Range("I38").Value = 2

# --- Row 39 ---
$ws2.Range("B39").Value = 2
$ws2.Range("C39").Value = "Eric Gasull"
$ws2.Range("D39").Value = "ericgasullserrano"
$ws2.Range("E39").Value = "Vista exposicions"
$ws2.Range("F39").Value = "Leonard Craciun"
$ws2.Range("G39").Value = 2
$ws2.Range("H39").Value = 1
$ws2.Range("I39").Value = 3

# --- Row 40 ---
$ws2.Range("B40").Value = 2
$ws2.Range("C40").Value = "Eric Gasull"
$ws2.Range("D40").Value = "ericgasullserrano"
$ws2.Range("E40").Value = "Vista crear exposicions"
$ws2.Range("F40").Value = "Leonard Craciun"
$ws2.Range("G40").Value = 2.5
$ws2.Range("H40").Value = 1
$ws2.Range("I40").Value = 2.5

# --- Row 41 ---
$ws2.Range("B41").Value = 2
$ws2.Range("C41").Value = "Eric Gasull"
$ws2.Range("D41").Value = "ericgasullserrano"
$ws2.Range("E41").Value = "Logica crear exposicions"
$ws2.Range("F41").Value = "Leonard Craciun"
$ws2.Range("G41").Value = 2
$ws2.Range("H41").Value = 1
$ws2.Range("I41").Value = 2

# --- Update the saved selection/active cell on Sprint2 ---
$ws2.Activate()
$ws2.Range("E26").Select()
